# lec2-1-function-templates.pptx
#
# The commit ("Adding iterator lecture project and massive update of
# others") touches this deck only cosmetically: the two C++ code snippets
# on the "demonstration of instantiation" slide get switched to a
# monospaced "Consolas" font. (The rest of the upstream diff - xmlns
# attribute reordering inside the hiddenFill/hiddenLine extLst boilerplate
# and the datetimeFigureOut placeholder text on every slideLayout/master -
# is simply what real PowerPoint re-emits whenever it resaves the whole
# package; it carries no addressable content and isn't reachable through
# the object model, so it's left alone here.)

$p = $ppt.ActivePresentation

# "демонстрация инстанцирования" slide.
$slide = $p.Slides.Item(46)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 2: extern template int max<int> (int, int);
$tr.Paragraphs(2, 1).Font.Name = "Consolas"

# Paragraph 4: template int max<int>(int, int);
$tr.Paragraphs(4, 1).Font.Name = "Consolas"
